$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.097369194030762
$ws.Range("B1").Value = 4.433632373809814
$ws.Range("C1").Value = 4.509809017181396
$ws.Range("D1").Value = 1.559965491294861
$ws.Range("E1").Value = 1.42816162109375
